$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "75.893.86"
$ws.Range("E2").Value = "  +1.44%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.908.00"
$ws.Range("E3").Value = "  +3.35%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "202.12"
$ws.Range("E5").Value = "  +7.72%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "597.10"
$ws.Range("E6").Value = "  +0.43%  "
$ws.Range("E7").Value = "  -0.06%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("E9").Value = "  +3.31%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "2.909.02"
$ws.Range("E10").Value = "  +3.39%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.431"
$ws.Range("E11").Value = "  +16.46%  "
$ws.Range("E12").Value = "  +0.44%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.88"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.447.32"
$ws.Range("E14").Value = "  +3.24%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "75.852.01"
$ws.Range("E15").Value = "  +1.36%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.70"
$ws.Range("E16").Value = "  +3.40%  "
$ws.Range("E17").Value = "  +1.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.909.64"
$ws.Range("E18").Value = "  +3.11%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.93"
$ws.Range("E19").Value = "  +5.42%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.72"
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "371.98"
$ws.Range("E21").Value = "  -0.89%  "
$ws.Range("E22").Value = "  +3.59%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.31"
$ws.Range("E23").Value = "  +6.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.00"
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "70.97"
$ws.Range("E25").Value = "  +0.47%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.034.51"
$ws.Range("E26").Value = "  +2.76%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "4.22"
$ws.Range("E27").Value = "  +1.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.67"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000108"
$ws.Range("E29").Value = "  +3.98%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.998"
$ws.Range("E30").Value = "  -0.14%  "
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "499.91"
$ws.Range("E32").Value = "  -2.47%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.67"
$ws.Range("E33").Value = "  -0.25%  "
$ws.Range("E34").Value = "  +1.78%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "166.43"
$ws.Range("E36").Value = "  +2.79%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.13"
$ws.Range("E37").Value = "  +1.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.61"
$ws.Range("E38").Value = "  +1.10%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.104"
$ws.Range("E39").Value = "  +22.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.112"
$ws.Range("E40").Value = "  -5.48%  "
$ws.Range("E41").Value = "  -0.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "180.59"
$ws.Range("E42").Value = "  -2.47%  "
$ws.Range("E43").Value = "  +2.67%  "
$ws.Range("E44").Value = "  -0.55%  "
$ws.Range("E45").Value = "  -0.68%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "40.16"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("E47").Value = "  -0.88%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.35"
$ws.Range("E48").Value = "  +1.14%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.572"
$ws.Range("E49").Value = "  +0.32%  "
$ws.Range("E50").Value = "  +0.48%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.654"
$ws.Range("E51").Value = "  +3.22%  "
